$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 0.3841813949771051
$arr[0,1] = 0.2733192036272456
$arr[0,2] = 0.0392471790369413
$arr[0,3] = 0
$arr[0,4] = 0.6906807884823252
$arr[0,5] = 0.5326446118630841
$arr[0,6] = 0.6826976885089806
$arr[0,7] = 0.7357532878019271
$arr[0,8] = 0
$arr[0,9] = 0.2678814174477679
$arr[0,10] = 0.3110847836143478
$arr[0,11] = 0.1510427188042094
$arr[0,12] = 0
$arr[0,13] = 2.400584438531723
$arr[1,0] = 0.3458251020513217
$arr[1,1] = 0.2738270911889416
$arr[1,2] = 0.03610526849424645
$arr[1,3] = 0
$arr[1,4] = 0.6931245927083722
$arr[1,5] = 0.536538401606748
$arr[1,6] = 0.6879135275748283
$arr[1,7] = 0.7430658348067816
$arr[1,8] = 0
$arr[1,9] = 0.2343053385667133
$arr[1,10] = 0.3079669980419126
$arr[1,11] = 0.1430830947827531
$arr[1,12] = 0
$arr[1,13] = 2.41953041521424
$arr[2,0] = 0.3222815573456046
$arr[2,1] = 0.2741813335233729
$arr[2,2] = 0.0341599196603255
$arr[2,3] = 0
$arr[2,4] = 0.6950288778246829
$arr[2,5] = 0.5392695018689437
$arr[2,6] = 0.6913878795751955
$arr[2,7] = 0.747913161965208
$arr[2,8] = 0
$arr[2,9] = 0.2135959185306149
$arr[2,10] = 0.3062136218520664
$arr[2,11] = 0.1382456080987176
$arr[2,12] = 0
$arr[2,13] = 2.432445811456205
$arr[3,0] = 0.3126899281033673
$arr[3,1] = 0.2743363957167091
$arr[3,2] = 0.03336313967248117
$arr[3,3] = 0
$arr[3,4] = 0.695906488880496
$arr[3,5] = 0.5404680061278313
$arr[3,6] = 0.692872121064859
$arr[3,7] = 0.7499783879513053
$arr[3,8] = 0
$arr[3,9] = 0.2051336972745759
$arr[3,10] = 0.3055396723092869
$arr[3,11] = 0.1362869541615161
$arr[3,12] = 0
$arr[3,13] = 2.438031542693054
$arr[4,0] = 0.3110974204585375
$arr[4,1] = 0.2743627915767419
$arr[4,2] = 0.03323059245759197
$arr[4,3] = 0
$arr[4,4] = 0.6960583538077003
$arr[4,5] = 0.5406721840301785
$arr[4,6] = 0.6931227122938992
$arr[4,7] = 0.7503267478584021
$arr[4,8] = 0
$arr[4,9] = 0.2037271815389801
$arr[4,10] = 0.3054302162484746
$arr[4,11] = 0.1359624902843528
$arr[4,12] = 0
$arr[4,13] = 2.438978536723283
$arr[5,0] = 0.3221521899961033
$arr[5,1] = 0.2741833813412811
$arr[5,2] = 0.03414919028976016
$arr[5,3] = 0
$arr[5,4] = 0.6950403021348919
$arr[5,5] = 0.53928531890827
$arr[5,6] = 0.6914076194770828
$arr[5,7] = 0.7479406502784585
$arr[5,8] = 0
$arr[5,9] = 0.2134818862880792
$arr[5,10] = 0.3062043683629412
$arr[5,11] = 0.1382191415431642
$arr[5,12] = 0
$arr[5,13] = 2.432519836230384
$arr[6,0] = 0.3709550235922734
$arr[6,1] = 0.2734855501616025
$arr[6,2] = 0.03816723427830482
$arr[6,3] = 0
$arr[6,4] = 0.691439646106879
$arr[6,5] = 0.5339165542776954
$arr[6,6] = 0.684439736874161
$arr[6,7] = 0.7382004960451773
$arr[6,8] = 0
$arr[6,9] = 0.2563241349194527
$arr[6,10] = 0.3099764178443962
$arr[6,11] = 0.1482880011207541
$arr[6,12] = 0
$arr[6,13] = 2.406850895758041
$arr[7,0] = 0.4666895690888566
$arr[7,1] = 0.2724516023458179
$arr[7,2] = 0.04591662803539265
$arr[7,3] = 0
$arr[7,4] = 0.68758031429536
$arr[7,5] = 0.526089376998506
$arr[7,6] = 0.6729294567964388
$arr[7,7] = 0.7219344249526252
$arr[7,8] = 0
$arr[7,9] = 0.3395742586411927
$arr[7,10] = 0.318647228780037
$arr[7,11] = 0.1684222683381122
$arr[7,12] = 0
$arr[7,13] = 2.366686045423236
$arr[8,0] = 0.5370171503699055
$arr[8,1] = 0.2718933808169766
$arr[8,2] = 0.05152946395752167
$arr[8,3] = 0
$arr[8,4] = 0.6866944910370094
$arr[8,5] = 0.521986940465446
$arr[8,6] = 0.6657818757092855
$arr[8,7] = 0.7117096645511864
$arr[8,8] = 0
$arr[8,9] = 0.4002502925899307
$arr[8,10] = 0.3257909091854287
$arr[8,11] = 0.1834463842349123
$arr[8,12] = 0
$arr[8,13] = 2.343373488960793
$arr[9,0] = 0.5690035443106751
$arr[9,1] = 0.2716826476585723
$arr[9,2] = 0.05406510643202012
$arr[9,3] = 0
$arr[9,4] = 0.6867144539479995
$arr[9,5] = 0.5204788975432706
$arr[9,6] = 0.6628136725102749
$arr[9,7] = 0.707432512644413
$arr[9,8] = 0
$arr[9,9] = 0.4277430590670406
$arr[9,10] = 0.3292080121311471
$arr[9,11] = 0.1903303609719131
$arr[9,12] = 0
$arr[9,13] = 2.334112662912304
$arr[10,0] = 0.5811144772930845
$arr[10,1] = 0.2716090182268331
$arr[10,2] = 0.05502271324819219
$arr[10,3] = 0
$arr[10,4] = 0.6867827792101551
$arr[10,5] = 0.5199593696078111
$arr[10,6] = 0.6617303610301448
$arr[10,7] = 0.705866641098023
$arr[10,8] = 0
$arr[10,9] = 0.4381376723429469
$arr[10,10] = 0.3305259751929128
$arr[10,11] = 0.1929441201496473
$arr[10,12] = 0
$arr[10,13] = 2.330799041320205
$arr[11,0] = 0.5785062553435694
$arr[11,1] = 0.271624601867309
$arr[11,2] = 0.0548165912276346
$arr[11,3] = 0
$arr[11,4] = 0.6867653623919665
$arr[11,5] = 0.5200689667366944
$arr[11,6] = 0.6619618629769022
$arr[11,7] = 0.706201487225421
$arr[11,8] = 0
$arr[11,9] = 0.4358997411803784
$arr[11,10] = 0.3302410632137622
$arr[11,11] = 0.1923808939659608
$arr[11,12] = 0
$arr[11,13] = 2.331504094700335
$arr[12,0] = 0.56999995488124
$arr[12,1] = 0.2716764666922487
$arr[12,2] = 0.05414394141030243
$arr[12,3] = 0
$arr[12,4] = 0.6867188575677545
$arr[12,5] = 0.5204351226574602
$arr[12,6] = 0.662723732772335
$arr[12,7] = 0.7073026094046604
$arr[12,8] = 0
$arr[12,9] = 0.4285985603774805
$arr[12,10] = 0.3293159616759738
$arr[12,11] = 0.1905452582180303
$arr[12,12] = 0
$arr[12,13] = 2.333836176431447
$arr[13,0] = 0.5647893666728692
$arr[13,1] = 0.2717090377491758
$arr[13,2] = 0.05373158584581006
$arr[13,3] = 0
$arr[13,4] = 0.6866982839847466
$arr[13,5] = 0.5206661160062254
$arr[13,6] = 0.6631956965965387
$arr[13,7] = 0.7079840837298264
$arr[13,8] = 0
$arr[13,9] = 0.4241242351926644
$arr[13,10] = 0.3287524304635951
$arr[13,11] = 0.1894217781752303
$arr[13,12] = 0
$arr[13,13] = 2.335289809814242
$arr[14,0] = 0.534926578980162
$arr[14,1] = 0.2719080185021454
$arr[14,2] = 0.05136339438710991
$arr[14,3] = 0
$arr[14,4] = 0.686701692493898
$arr[14,5] = 0.5220927046128097
$arr[14,6] = 0.6659815504819022
$arr[14,7] = 0.7119967144986532
$arr[14,8] = 0
$arr[14,9] = 0.3984513281669138
$arr[14,10] = 0.3255709542250287
$arr[14,11] = 0.1829974810375461
$arr[14,12] = 0
$arr[14,13] = 2.344005750693114
$arr[15,0] = 0.5166046536638476
$arr[15,1] = 0.2720411223936878
$arr[15,2] = 0.04990602872021555
$arr[15,3] = 0
$arr[15,4] = 0.6868120697912445
$arr[15,5] = 0.5230596295700423
$arr[15,6] = 0.6677630969367243
$arr[15,7] = 0.7145541546704095
$arr[15,8] = 0
$arr[15,9] = 0.3826734572633939
$arr[15,10] = 0.3236620328178361
$arr[15,11] = 0.1790689283888156
$arr[15,12] = 0
$arr[15,13] = 2.349696946406482
$arr[16,0] = 0.5060658497737052
$arr[16,1] = 0.272121749352209
$arr[16,2] = 0.0490661305888267
$arr[16,3] = 0
$arr[16,4] = 0.6869153670243904
$arr[16,5] = 0.5236494907256954
$arr[16,6] = 0.668814460825196
$arr[16,7] = 0.716060343218146
$arr[16,8] = 0
$arr[16,9] = 0.3735882142415505
$arr[16,10] = 0.3225798344110586
$arr[16,11] = 0.1768139911570898
$arr[16,12] = 0
$arr[16,13] = 2.35309689395109
$arr[17,0] = 0.5024975262632836
$arr[17,1] = 0.272149748418812
$arr[17,2] = 0.04878147180981784
$arr[17,3] = 0
$arr[17,4] = 0.6869571808234127
$arr[17,5] = 0.5238549966330055
$arr[17,6] = 0.6691750159479142
$arr[17,7] = 0.7165763615145906
$arr[17,8] = 0
$arr[17,9] = 0.3705103690543865
$arr[17,10] = 0.3222161306956224
$arr[17,11] = 0.1760513138389328
$arr[17,12] = 0
$arr[17,13] = 2.354269789233911
$arr[18,0] = 0.5185551129202111
$arr[18,1] = 0.2720265324252935
$arr[18,2] = 0.05006133991366823
$arr[18,3] = 0
$arr[18,4] = 0.6867962000017371
$arr[18,5] = 0.5229532095250349
$arr[18,6] = 0.6675706888662631
$arr[18,7] = 0.7142782660710232
$arr[18,8] = 0
$arr[18,9] = 0.3843541016932761
$arr[18,10] = 0.3238636100499974
$arr[18,11] = 0.1794866481596671
$arr[18,12] = 0
$arr[18,13] = 2.349078014594866
$arr[19,0] = 0.5724985115525953
$arr[19,1] = 0.2716610656001919
$arr[19,2] = 0.05434158551412338
$arr[19,3] = 0
$arr[19,4] = 0.686730868408155
$arr[19,5] = 0.5203261749038433
$arr[19,6] = 0.6624988494344919
$arr[19,7] = 0.7069777232676842
$arr[19,8] = 0
$arr[19,9] = 0.4307435414702638
$arr[19,10] = 0.3295870363432982
$arr[19,11] = 0.1910842416929981
$arr[19,12] = 0
$arr[19,13] = 2.333145942999025
$arr[20,0] = 0.6077438449165413
$arr[20,1] = 0.2714581598410319
$arr[20,2] = 0.05712387918507034
$arr[20,3] = 0
$arr[20,4] = 0.6870423352611468
$arr[20,5] = 0.5189096556575095
$arr[20,6] = 0.6594212153138059
$arr[20,7] = 0.7025199396501307
$arr[20,8] = 0
$arr[20,9] = 0.4609664236498645
$arr[20,10] = 0.3334673587184795
$arr[20,11] = 0.1987043545280329
$arr[20,12] = 0
$arr[20,13] = 2.32385981624229
$arr[21,0] = 0.5889338945977158
$arr[21,1] = 0.2715631790482504
$arr[21,2] = 0.05564031297160454
$arr[21,3] = 0
$arr[21,4] = 0.6868437105327629
$arr[21,5] = 0.5196381824018275
$arr[21,6] = 0.6610421280989272
$arr[21,7] = 0.704870456213559
$arr[21,8] = 0
$arr[21,9] = 0.4448448259790894
$arr[21,10] = 0.3313836028308259
$arr[21,11] = 0.1946337123100079
$arr[21,12] = 0
$arr[21,13] = 2.328712943967417
$arr[22,0] = 0.5176733266955296
$arr[22,1] = 0.2720331157653817
$arr[22,2] = 0.04999113007148992
$arr[22,3] = 0
$arr[22,4] = 0.6868032506326642
$arr[22,5] = 0.5230012162374962
$arr[22,6] = 0.6676575920459342
$arr[22,7] = 0.7144028835364367
$arr[22,8] = 0
$arr[22,9] = 0.3835943269513393
$arr[22,10] = 0.3237724294253042
$arr[22,11] = 0.1792977857011948
$arr[22,12] = 0
$arr[22,13] = 2.349357434944395
$arr[23,0] = 0.4407903076588866
$arr[23,1] = 0.2726957620866912
$arr[23,2] = 0.04383427172848542
$arr[23,3] = 0
$arr[23,4] = 0.688281833909258
$arr[23,5] = 0.5279175101256257
$arr[23,6] = 0.6758131305989252
$arr[23,7] = 0.726031595646468
$arr[23,8] = 0
$arr[23,9] = 0.3171369478827444
$arr[23,10] = 0.3161655144081976
$arr[23,11] = 0.1629343454752714
$arr[23,12] = 0
$arr[23,13] = 2.376463192267977
$ws.Range("B2:O25").Value = $arr
Write-Host "done"
